# Update odds values in Sheet1 to reflect latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 2 (Corinthians vs Athletico-PR)
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3.2
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8

# Row 3 (Flamengo RJ vs Fluminense)
$ws.Range("V3").Value = 1.67

# Row 4 (America MG vs Goias)
$ws.Range("V4").Value = 1.67

# Row 8 (Union Comercio vs Comerciantes Unidos)
$ws.Range("J8").Value = 2.62
$ws.Range("Q8").Value = 1.79
$ws.Range("R8").Value = 1.94

# Row 10 (Sport Huancayo vs Grau)
$ws.Range("G10").Value = 1.86
$ws.Range("J10").Value = 2.62
$ws.Range("N10").Value = 9

# Row 12 (Fenix vs CA Cerro)
$ws.Range("K12").Value = 1.92

# Row 13 (Nacional vs Miramar)
$ws.Range("Q13").Value = 1.84
$ws.Range("R13").Value = 1.89
